$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.693.75'
$ws.Range('E2').Value = '  -0.59%  '
$ws.Range('D3').Value = '3.618.34'
$ws.Range('E3').Value = '  -1.20%  '
$ws.Range('D4').Value = '''0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').Value = '''585.96'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.61%  '
$ws.Range('D6').Value = '''195.03'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.17%  '
$ws.Range('D7').Value = '3.615.11'
$ws.Range('E7').Value = '  -1.16%  '
$ws.Range('D8').Value = '''0.620'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').Value = '''0.683'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.05%  '
$ws.Range('D11').Value = '''0.152'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('D12').Value = '''55.34'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.67%  '
$ws.Range('D13').Value = '''0.0000278'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.51%  '
$ws.Range('D14').Value = '''10.04'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.66%  '
$ws.Range('D15').Value = '4.180.13'
$ws.Range('E15').Value = '  -1.63%  '
$ws.Range('D16').Value = '3.612.62'
$ws.Range('E16').Value = '  -1.33%  '
$ws.Range('D17').Value = '''0.126'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.62%  '
$ws.Range('D18').Value = '''12.43'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.28%  '
$ws.Range('D19').Value = '''18.55'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('D20').Value = '67.625.28'
$ws.Range('E20').Value = '  -0.39%  '
$ws.Range('E21').Value = '  -2.32%  '
$ws.Range('D22').Value = '''404.57'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.05%  '
$ws.Range('D23').Value = '''13.35'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +20.15%  '
$ws.Range('D24').Value = '''4.27'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.59%  '
$ws.Range('D25').Value = '''86.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.66%  '
$ws.Range('D26').Value = '''2.95'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').Value = '''12.68'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = '''3.85'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.90%  '
$ws.Range('B29').Value = 'LEO'
$ws.Range('C29').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D29').Value = '''6.11'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.81%  '
$ws.Range('D30').Value = '''8.12'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +12.11%  '
$ws.Range('D31').Value = '''9.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.99%  '
$ws.Range('D32').Value = '''31.52'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.08%  '
$ws.Range('D33').Value = '''669.72'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +10.16%  '
$ws.Range('D34').Value = '''12.30'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').Value = '''0.118'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.47%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').Value = '''43.21'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.32%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = '''64.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.82%  '
$ws.Range('D38').Value = '''0.426'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +8.28%  '
$ws.Range('E39').Value = '  +0.36%  '
$ws.Range('D40').Value = '0.0₃0798'
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('D41').Value = '''2.97'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +17.51%  '
$ws.Range('E42').Value = '  +7.93%  '
$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '3.167.74'
$ws.Range('E43').Value = '  +14.35%  '
$ws.Range('B44').Value = 'Kaspa'
$ws.Range('C44').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D44').Value = '''0.134'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.82%  '
$ws.Range('D45').Value = '''0.997'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.43%  '
$ws.Range('D46').Value = '''0.0421'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.00%  '
$ws.Range('B47').Value = 'Stellar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D47').Value = '''0.132'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.46%  '
$ws.Range('B48').Value = 'THORChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D48').Value = '''8.84'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.31%  '
$ws.Range('B49').Value = 'ApeXProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D49').Value = '''3.10'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('D50').Value = '''142.92'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''2.74'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.95%  '
